# Updated UTR OTD input files
# Apply value changes to the "Extent of Contamination" sheet (Table1568),
# plus the selected-cell change on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extent of Contamination")

# Row 4 (Indoor / Area Contaminated): Parameter 1 (G) 0 -> 747.8695
$ws.Range("G4").Value = 747.86950000000002

# Row 5 (Indoor / Loading): Distribution Type Constant -> Uniform,
# Parameter 1 -1 -> 3.0319702616531905, Parameter 2 (blank) -> 8.0319702616531909
$ws.Range("F5").Value = "Uniform"
$ws.Range("G5").Value = 3.0319702616531905
$ws.Range("H5").Value = 8.0319702616531909

# Row 6 (Underground / Area Contaminated): Parameter 1 2682.854 -> 0
$ws.Range("G6").Value = 0

# Row 7 (Underground / Loading): Distribution Type Uniform -> Constant,
# Parameter 1 5.6318533337268804 -> -1, Parameter 2 6.1323408067707499 -> blank
# G7 previously had no explicit cell style (it held the Uniform min); bring
# it into line with the rest of the Parameter-1 column (style like G6) now
# that it holds a normal Constant parameter value.
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("F7").Value = "Constant"
$ws.Range("G7").Value = -1
$ws.Range("H7").ClearContents()

# Row 9 (Indoor / Commercial breakout): Parameter 1 0 -> 1
$ws.Range("G9").Value = 1

# Underground surface-type breakout rows 21-26
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 0

# Indoor surface-type breakout rows 27-32
$ws.Range("G27").Value = 0.5
$ws.Range("G28").Value = 0.125
$ws.Range("G29").Value = 0.0625
$ws.Range("G30").Value = 0.0625
$ws.Range("G31").Value = 0.125
$ws.Range("G32").Value = 0.125

# The conditional formatting that covered A2:L6 A8:L32 A7:F7 H7:L7 (i.e. all
# rows except G7/H7, which used to hold the Uniform-distribution parameters)
# now applies to the full, contiguous A2:L32 block.
$fc = $ws.Range("A2:L6").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:L32"))

# Update the selected cell shown when the sheet is active
$ws.Range("G7").Select()
